$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before column B (shifts Spezialpreis/Artikelname/Verkaufspreis right) ---
$ws.Columns.Item(2).Insert()

# --- 2. Grow the table to cover the freshly inserted column BEFORE touching the
#        header cells, so the table's column-name cache resyncs from the cells
#        we are about to (re)write instead of keeping stale names/appending an
#        extra "Column5".
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E3"))

# --- 3. Header row: restore the header text that moved right, then set the new header ---
$ws.Range("C1").Value = "Spezialpreis"
$ws.Range("D1").Value = "Artikelname"
$ws.Range("E1").Value = "Verkaufspreis"

# --- 4. Fill in the new "Suisanummer" column ---
# Row 2 keeps a genuine number (1019.343): write the value BEFORE applying the
# text number format so it stays a numeric cell, just displayed via style 3.
$ws.Range("B2").Value = 1019.343

# Row 3 must be stored as literal text ("1021.174"): apply the text format
# FIRST so the numeric-looking string isn't auto-converted to a number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1021.174"

# Header + row2 pick up the same text format afterwards.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "Suisanummer"
$ws.Range("B2").NumberFormat = "@"

# Match column B's width to column A and mark the whole column with the text style.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth()

# --- 5. Move the active selection the way the author left it ---
$ws.Range("B5").Select()
